$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 116
$ws1.Range("F5").Value = 5089
$ws1.Range("F7").Value = 628
$ws1.Range("F8").Value = 298
$ws1.Range("F9").Value = 766
$ws1.Range("F10").Value = 251

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 32
$ws2.Range("F3").Value = 4

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 116
$ws4.Range("F5").Value = 5089
$ws4.Range("F7").Value = 628
$ws4.Range("F8").Value = 298
$ws4.Range("F9").Value = 766
$ws4.Range("F10").Value = 32
$ws4.Range("F11").Value = 251
$ws4.Range("F13").Value = 4
